# Datentreiber.xlsx - add "FreelancerMap" / "Testmanager" as a new first
# data row on the "Testdaten" sheet (config.xml / Test.java / VakanzenGrabber
# now recognise the FreelancerMap portal), keep "FreelanceDE" / "Testanalyst"
# as the following row, and drop the two trailing blank filler rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testdaten")

# Push the existing data row (FreelanceDE / Testanalyst) down into row 3 -
# row 3 already carries the styled ("s=1") column-A formatting in the
# original sheet, so writing the values directly here (instead of using
# Rows.Insert, which does not propagate that formatting in this runtime)
# keeps the look identical to the source file.
$ws.Range("A3").Value = "FreelanceDE"
$ws.Range("B3").Value = "Testanalyst"

# Row 2 becomes the new FreelancerMap entry.
$ws.Range("A2").Value = "FreelancerMap"
$ws.Range("B2").Value = "Testmanager"

# The two trailing empty filler rows (old rows 3 & 4) collapse to one row,
# which is now fully populated above - remove the now-superfluous row 4.
$ws.Rows(4).Delete()

# Match the saved selection / view state from the edit.
$ws.Range("D18").Select()

# The saved file also records an explicit portrait page setup.
$ws.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait
